$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add city1..city6 values in column E ---------------------------------
# E1 already holds the shared string that used to read "Music"; giving it
# the value "city1" reuses that same shared-string slot (keeping its
# existing style), matching the source edit which simply renamed that
# string in place.
$ws.Range("E1").Value = "city1"
$ws.Range("E2").Value = "city2"
$ws.Range("E3").Value = "city3"
$ws.Range("E4").Value = "city4"
$ws.Range("E5").Value = "city5"
$ws.Range("E6").Value = "city6"

# Give the new E2:E6 cells the same formatting used by the rest of the
# table (column D), i.e. the Calibri 10 / theme text color used throughout
# the sheet.
$ws.Range("D2:D6").Copy()
$ws.Range("E2:E6").PasteSpecial(-4122)

# --- Register the extra font used for phonetic guides ---------------------
# The workbook gains a 3rd font (Arial 8) that is referenced by the new
# phoneticPr on the sheet. Create it through a throw-away cell style so the
# font is registered in the workbook without altering any cell's style.
$fontStyle = $wb.Styles.Add("__TempPhoneticFont")
$fontStyle.Font.Name = "Arial"
$fontStyle.Font.Size = 8
$fontStyle.Delete()

# --- Update the selected cell shown when the sheet is opened --------------
$ws.Range("I7").Select()
